$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data in columns A and B (A had 0..6 index, B had the
# "Year" header + the actual year values). We rebuild the table in column A
# only, as the per-capita-income prediction/CSV work removed the helper
# index column and moved "Year" + values into column A.
$ws.Range("A1:B8").ClearContents()

$ws.Range("A1").Value = "Year"
$ws.Range("A2").Value = 2024
$ws.Range("A3").Value = 2025
$ws.Range("A4").Value = 2026
$ws.Range("A5").Value = 2027
$ws.Range("A6").Value = 2028
$ws.Range("A7").Value = 2029
$ws.Range("A8").Value = 2030

$ws.Range("N12").Select()
